# The source workbook tracks one new weekly price observation for Perejil
# (Parsley) at Vega Central Mapocho de Santiago. The new record is inserted
# as row 426 (sorted position within the date series), which pushes the
# former rows 426-501 down to 427-502.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 426; this shifts existing rows 426..501 down
# to 427..502 and copies formatting from the row above (matches how the
# D column keeps its date style).
$ws.Rows.Item(426).Insert()

# Populate the newly inserted row 426 with the new observation.
$ws.Range("A426").Value = 9
$ws.Range("B426").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C426").Value = "Metropolitana"
$ws.Range("D426").Value = 44951
$ws.Range("E426").Value = 13
$ws.Range("F426").Value = 100112044
$ws.Range("G426").Value = "Perejil"
$ws.Range("H426").Value = "Sin especificar"
$ws.Range("I426").Value = "Primera"
$ws.Range("J426").Value = 70
$ws.Range("K426").Value = 14000
$ws.Range("L426").Value = 16000
$ws.Range("M426").Value = 15000
$ws.Range("N426").Value = "$/docena de atados"
$ws.Range("O426").Value = "Región Metropolitana"
$ws.Range("P426").Value = 5000
$ws.Range("Q426").Value = 3
$ws.Range("R426").Value = "Hortaliza"
